$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values for rows 2-8 (columns B, C, D); column A (iteration index) is unchanged
$ws.Range("B2").Value = -2
$ws.Range("C2").Value = -0.909297426825682
$ws.Range("D2").Value = 1.0005

$ws.Range("B3").Value = 3
$ws.Range("C3").Value = 0.141120008059867
$ws.Range("D3").Value = 1.0005

$ws.Range("B4").Value = 2.32826701379322
$ws.Range("C4").Value = 0.726576188310738
$ws.Range("D4").Value = 0.671732986206781

$ws.Range("B5").Value = 3.16191641257756
$ws.Range("C5").Value = -0.0203223598776722
$ws.Range("D5").Value = 0.833649398784336

$ws.Range("B6").Value = 3.13923364912194
$ws.Range("C6").Value = 0.002359002279916
$ws.Range("D6").Value = 0.0226827634556157

$ws.Range("B7").Value = 3.14159279714554
$ws.Range("C7").Value = -0.000000143555744946614
$ws.Range("D7").Value = 0.0023591480235984

$ws.Range("B8").Value = 3.14159265358966
$ws.Range("C8").Value = 0.000000000000133349227634934
$ws.Range("D8").Value = 0.000000143555878295842

# Rows 9 and 10 no longer belong to the table - remove them entirely
$ws.Range("A9:D10").Delete()
